$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 38760
$ws.Range("J126").Value = 38760
$ws.Range("L126").Value = 38760
$ws.Range("N126").Value = -48640
$ws.Range("H130").Value = 97850
$ws.Range("J130").Value = 97850
$ws.Range("L130").Value = 97850
$ws.Range("N130").Value = -107890
$ws.Range("H141").Value = 6364
$ws.Range("I141").Value = 3682.5
$ws.Range("J141").Value = 9428.571
$ws.Range("K141").Value = 11047.5
$ws.Range("L141").Value = 28285.713
$ws.Range("M141").Value = -5867.5
$ws.Range("N141").Value = -38645.713

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4192.6875
$ws.Range("I2").Value = 4807.5454
$ws.Range("J2").Value = 2840
$ws.Range("K2").Value = 4807.5454
$ws.Range("L2").Value = 2840
$ws.Range("M2").Value = -4694.5454
$ws.Range("N2").Value = -3066
$ws.Range("H45").Value = 13927.875
$ws.Range("I45").Value = 25831
$ws.Range("J45").Value = 2024.75
$ws.Range("K45").Value = 25831
$ws.Range("L45").Value = 2024.75
$ws.Range("M45").Value = -25454
$ws.Range("N45").Value = -2778.75
$ws.Range("H116").Value = 4192.6875
$ws.Range("I116").Value = 4807.5454
$ws.Range("J116").Value = 2840
$ws.Range("K116").Value = 4807.5454
$ws.Range("L116").Value = 2840
$ws.Range("M116").Value = -2513.5454
$ws.Range("N116").Value = -7428

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4192.6875
$ws.Range("I3").Value = 4807.5454
$ws.Range("J3").Value = 2840
$ws.Range("K3").Value = 4807.5454
$ws.Range("L3").Value = 2840
$ws.Range("M3").Value = -4693.5454
$ws.Range("N3").Value = -3068

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9437290
$ws.Range("I31").Value = 22728470
$ws.Range("K31").Value = 22728470
$ws.Range("M31").Value = -22728175
$ws.Range("H34").Value = 9437290
$ws.Range("I34").Value = 22728470
$ws.Range("K34").Value = 22728470
$ws.Range("M34").Value = -22728268

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11610.223
$ws.Range("I68").Value = 20388.4
$ws.Range("J68").Value = 637.5
$ws.Range("K68").Value = 61165.2
$ws.Range("L68").Value = 1912.5
$ws.Range("M68").Value = -60354.2
$ws.Range("N68").Value = -3534.5
$ws.Range("H70").Value = 3576.5881
$ws.Range("I70").Value = 1543.1428
$ws.Range("K70").Value = 4629.428400000001
$ws.Range("M70").Value = -4314.428400000001
$ws.Range("H71").Value = 11610.223
$ws.Range("I71").Value = 20388.4
$ws.Range("J71").Value = 637.5
$ws.Range("K71").Value = 183495.6
$ws.Range("L71").Value = 5737.5
$ws.Range("M71").Value = -179439.6
$ws.Range("N71").Value = -13849.5
$ws.Range("H73").Value = 3576.5881
$ws.Range("I73").Value = 1543.1428
$ws.Range("K73").Value = 4629.428400000001
$ws.Range("M73").Value = -3537.428400000001
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = $null
$ws.Range("N75").Value = -4996
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 9000
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = -18984
$ws.Range("H113").Value = 465.77274
$ws.Range("I113").Value = 439.7
$ws.Range("J113").Value = 487.5
$ws.Range("K113").Value = 1319.1
$ws.Range("L113").Value = 1462.5
$ws.Range("M113").Value = 850.9000000000001
$ws.Range("N113").Value = -5802.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3335.261
$ws.Range("I80").Value = 2700.4167
$ws.Range("J80").Value = 4027.818
$ws.Range("K80").Value = 2700.4167
$ws.Range("L80").Value = 4027.818
$ws.Range("M80").Value = -1702.4167
$ws.Range("N80").Value = -6023.818
$ws.Range("H83").Value = 3335.261
$ws.Range("I83").Value = 2700.4167
$ws.Range("J83").Value = 4027.818
$ws.Range("K83").Value = 13502.0835
$ws.Range("L83").Value = 20139.09
$ws.Range("M83").Value = -8510.083500000001
$ws.Range("N83").Value = -30123.09
$ws.Range("H102").Value = 2247.077
$ws.Range("I102").Value = 1737.1333
$ws.Range("J102").Value = 2942.4546
$ws.Range("K102").Value = 1737.1333
$ws.Range("L102").Value = 2942.4546
$ws.Range("M102").Value = -115.1333
$ws.Range("N102").Value = -6186.4546

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 28253.166
$ws.Range("J130").Value = 24184.75
$ws.Range("L130").Value = 24184.75
$ws.Range("N130").Value = -34224.75
$ws.Range("H141").Value = 55947.6
$ws.Range("J141").Value = 55947.6
$ws.Range("L141").Value = 55947.6
$ws.Range("N141").Value = -66307.60000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4294285.5
$ws.Range("I2").Value = 50000000
$ws.Range("J2").Value = 778461.5600000001
$ws.Range("K2").Value = 50000000
$ws.Range("L2").Value = 778461.5600000001
$ws.Range("M2").Value = -49999888
$ws.Range("N2").Value = -778685.5600000001
$ws.Range("H3").Value = 6600
$ws.Range("J3").Value = 6600
$ws.Range("L3").Value = 6600
$ws.Range("N3").Value = -6828
$ws.Range("H4").Value = 229198.11
$ws.Range("J4").Value = 229198.11
$ws.Range("L4").Value = 229198.11
$ws.Range("N4").Value = -229424.11
$ws.Range("H5").Value = 1251875.2
$ws.Range("I5").Value = 500.5
$ws.Range("J5").Value = 2503250
$ws.Range("K5").Value = 500.5
$ws.Range("L5").Value = 2503250
$ws.Range("M5").Value = -388.5
$ws.Range("N5").Value = -2503474
$ws.Range("H6").Value = 18000.334
$ws.Range("I6").Value = 3001
$ws.Range("J6").Value = 25500
$ws.Range("K6").Value = 3001
$ws.Range("L6").Value = 25500
$ws.Range("M6").Value = -2886
$ws.Range("N6").Value = -25730
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = -3280
$ws.Range("H20").Value = 9383.125
$ws.Range("I20").Value = 6010
$ws.Range("J20").Value = 9865
$ws.Range("K20").Value = 6010
$ws.Range("L20").Value = 9865
$ws.Range("M20").Value = -5770
$ws.Range("N20").Value = -10345
